$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1916261221985991
$ws.Range("D2").Value = 0.1438557052619132
$ws.Range("E2").Value = 0.1523442560692985
$ws.Range("F2").Value = 1.856329772453911
$ws.Range("G2").Value = 0.002498313024655947
$ws.Range("I2").Value = 1.305044519253165
$ws.Range("J2").Value = 0.1958912783820281
$ws.Range("K2").Value = 1.651681354704237
$ws.Range("L2").Value = 0.2144890373025561
$ws.Range("N2").Value = 1.395055625168936
$ws.Range("O2").Value = 4.793580040525967
$ws.Range("C3").Value = 0.1885383805380485
$ws.Range("D3").Value = 0.1411154674472854
$ws.Range("E3").Value = 0.151341998071171
$ws.Range("F3").Value = 1.861361175236979
$ws.Range("G3").Value = 0.002501295677841409
$ws.Range("I3").Value = 1.309014802672166
$ws.Range("J3").Value = 0.1958234761985622
$ws.Range("K3").Value = 1.526098408042145
$ws.Range("L3").Value = 0.2140403703250158
$ws.Range("N3").Value = 1.398139420227572
$ws.Range("O3").Value = 4.817067947351916
$ws.Range("C4").Value = 0.1867104894596281
$ws.Range("D4").Value = 0.1394759377343036
$ws.Range("E4").Value = 0.1507813899653243
$ws.Range("F4").Value = 1.865406520259185
$ws.Range("G4").Value = 0.002503225797293504
$ws.Range("I4").Value = 1.312133890792786
$ws.Range("J4").Value = 0.1958584306702598
$ws.Range("K4").Value = 1.449186388637969
$ws.Range("L4").Value = 0.2138414519630487
$ws.Range("N4").Value = 1.400529146630561
$ws.Range("O4").Value = 4.834076669789226
$ws.Range("C5").Value = 0.1859827942360965
$ws.Range("D5").Value = 0.1388187044486102
$ws.Range("E5").Value = 0.1505667578261622
$ws.Range("F5").Value = 1.867295488322391
$ws.Range("G5").Value = 0.002504037244497925
$ws.Range("I5").Value = 1.313576287412673
$ws.Range("J5").Value = 0.1958919652562479
$ws.Range("K5").Value = 1.417895780953444
$ws.Range("L5").Value = 0.2137796904905116
$ws.Range("N5").Value = 1.401628083798343
$ws.Range("O5").Value = 4.841658248194307
$ws.Range("C6").Value = 0.1858630015194791
$ws.Range("D6").Value = 0.1387102313177877
$ws.Range("E6").Value = 0.1505319543483985
$ws.Range("F6").Value = 1.867623674620091
$ws.Range("G6").Value = 0.002504173491110604
$ws.Range("I6").Value = 1.313826145572975
$ws.Range("J6").Value = 0.1958986996936609
$ws.Range("K6").Value = 1.412703186858778
$ws.Range("L6").Value = 0.2137706020497987
$ws.Range("N6").Value = 1.401818126843651
$ws.Range("O6").Value = 4.842956440930607
$ws.Range("C7").Value = 0.18670060580952
$ws.Range("D7").Value = 0.139467029864683
$ws.Range("E7").Value = 0.1507784393440303
$ws.Range("F7").Value = 1.865431021874116
$ws.Range("G7").Value = 0.002503236639897943
$ws.Range("I7").Value = 1.31215264967237
$ws.Range("J7").Value = 0.1958588047802863
$ws.Range("K7").Value = 1.448764180355283
$ws.Range("L7").Value = 0.2138405408217352
$ws.Range("N7").Value = 1.400543460318723
$ws.Range("O7").Value = 4.834176284425098
$ws.Range("C8").Value = 0.1905474076199596
$ws.Range("D8").Value = 0.1429019954106394
$ws.Range("E8").Value = 0.1519873336475577
$ws.Range("F8").Value = 1.857866201005848
$ws.Range("G8").Value = 0.002499320990893071
$ws.Range("I8").Value = 1.306272074872346
$ws.Range("J8").Value = 0.1958520238628054
$ws.Range("K8").Value = 1.608340910501681
$ws.Range("L8").Value = 0.21431847317033
$ws.Range("N8").Value = 1.396016103053441
$ws.Range("O8").Value = 4.801141631342318
$ws.Range("C9").Value = 0.1986271790557055
$ws.Range("D9").Value = 0.1499760076709578
$ws.Range("E9").Value = 0.1547909312930322
$ws.Range("F9").Value = 1.850616842972045
$ws.Range("G9").Value = 0.002492422631049711
$ws.Range("I9").Value = 1.30014657005885
$ws.Range("J9").Value = 0.1964451687623736
$ws.Range("K9").Value = 1.92274486117276
$ws.Range("L9").Value = 0.2158614013415203
$ws.Range("N9").Value = 1.39106289031379
$ws.Range("O9").Value = 4.756898066690667
$ws.Range("C10").Value = 0.2048866222253309
$ws.Range("D10").Value = 0.1553759384977411
$ws.Range("E10").Value = 0.1571127464875275
$ws.Range("F10").Value = 1.849916775176283
$ws.Range("G10").Value = 0.002487825302122958
$ws.Range("I10").Value = 1.298944216692597
$ws.Range("J10").Value = 0.1972492926819314
$ws.Range("K10").Value = 2.154546799831337
$ws.Range("L10").Value = 0.2173621701655648
$ws.Range("N10").Value = 1.389801054595893
$ws.Range("O10").Value = 4.736931488222183
$ws.Range("C11").Value = 0.2078036697314332
$ws.Range("D11").Value = 0.1578758247103593
$ws.Range("E11").Value = 0.1582255092905207
$ws.Range("F11").Value = 1.850603398092915
$ws.Range("G11").Value = 0.002485835099233815
$ws.Range("I11").Value = 1.299113993922418
$ws.Range("J11").Value = 0.1976948236063194
$ws.Range("K11").Value = 2.260157770114176
$ws.Range("L11").Value = 0.2181242237902055
$ws.Range("N11").Value = 1.389740170013184
$ws.Range("O11").Value = 4.730575451298904
$ws.Range("C12").Value = 0.2089182103476901
$ws.Range("D12").Value = 0.1588286409093627
$ws.Range("E12").Value = 0.1586549769917411
$ws.Range("F12").Value = 1.851007942319541
$ws.Range("G12").Value = 0.002485095929553549
$ws.Range("I12").Value = 1.299281372230865
$ws.Range("J12").Value = 0.1978749724897995
$ws.Range("K12").Value = 2.300171369981115
$ws.Range("L12").Value = 0.2184241625643182
$ws.Range("N12").Value = 1.389790640656656
$ws.Range("O12").Value = 4.728561005863298
$ws.Range("C13").Value = 0.2086777344046169
$ws.Range("D13").Value = 0.1586231620665473
$ws.Range("E13").Value = 0.1585621244009658
$ws.Range("F13").Value = 1.850914388169414
$ws.Range("G13").Value = 0.002485254480244938
$ws.Range("I13").Value = 1.299240739128265
$ws.Range("J13").Value = 0.1978356660495351
$ws.Range("K13").Value = 2.291552829272234
$ws.Range("L13").Value = 0.2183590605982317
$ws.Range("N13").Value = 1.389776504883315
$ws.Range("O13").Value = 4.728977393263506
$ws.Range("C14").Value = 0.2078951653998047
$ws.Range("D14").Value = 0.1579540904133125
$ws.Range("E14").Value = 0.1582606800030071
$ws.Range("F14").Value = 1.850633783464843
$ws.Range("G14").Value = 0.002485773997543066
$ws.Range("I14").Value = 1.299125698086357
$ws.Range("J14").Value = 0.1977094155359538
$ws.Range("K14").Value = 2.26344930576613
$ws.Range("L14").Value = 0.2181486724031103
$ws.Range("N14").Value = 1.389742850158271
$ws.Range("O14").Value = 4.730401854536069
$ws.Range("C15").Value = 0.2074171082486629
$ws.Range("D15").Value = 0.1575450649608712
$ws.Range("E15").Value = 0.1580770886636529
$ws.Range("F15").Value = 1.850480727852414
$ws.Range("G15").Value = 0.002486094100320302
$ws.Range("I15").Value = 1.299068657796489
$ws.Range("J15").Value = 0.1976335718899591
$ws.Range("K15").Value = 2.246237751053627
$ws.Range("L15").Value = 0.2180212823840293
$ws.Range("N15").Value = 1.389731803329539
$ws.Range("O15").Value = 4.731325496255039
$ws.Range("C16").Value = 0.2046973790759807
$ws.Range("D16").Value = 0.1552134319979359
$ws.Range("E16").Value = 0.1570411591111167
$ws.Range("F16").Value = 1.849892128074742
$ws.Range("G16").Value = 0.002487957395839874
$ws.Range("I16").Value = 1.298947545735949
$ws.Range("J16").Value = 0.1972217780772993
$ws.Range("K16").Value = 2.147647948153235
$ws.Range("L16").Value = 0.2173139610300368
$ws.Range("N16").Value = 1.389815336554832
$ws.Range("O16").Value = 4.73740176964759
$ws.Range("C17").Value = 0.2030466744921569
$ws.Range("D17").Value = 0.1537941161005563
$ws.Range("E17").Value = 0.1564201040349928
$ws.Range("F17").Value = 1.849788471397531
$ws.Range("G17").Value = 0.002489126324611147
$ws.Range("I17").Value = 1.299056835659528
$ws.Range("J17").Value = 0.1969895560843682
$ws.Range("K17").Value = 2.087206336211977
$ws.Range("L17").Value = 0.2169003332242312
$ws.Range("N17").Value = 1.389997818454859
$ws.Range("O17").Value = 4.741828014339717
$ws.Range("C18").Value = 0.2021037888764852
$ws.Range("D18").Value = 0.1529818559019276
$ws.Range("E18").Value = 0.1560682156769708
$ws.Range("F18").Value = 1.849823454079853
$ws.Range("G18").Value = 0.002489808186074482
$ws.Range("I18").Value = 1.299187157381652
$ws.Range("J18").Value = 0.1968634945270082
$ws.Range("K18").Value = 2.052457417076596
$ws.Range("L18").Value = 0.2166698965451843
$ws.Range("N18").Value = 1.390151090679495
$ws.Range("O18").Value = 4.744630518418205
$ws.Range("C19").Value = 0.2017856727439522
$ws.Range("D19").Value = 0.1527075439856134
$ws.Range("E19").Value = 0.155949988439005
$ws.Range("F19").Value = 1.849851546682444
$ws.Range("G19").Value = 0.002490040690375236
$ws.Range("I19").Value = 1.299242867916881
$ws.Range("J19").Value = 0.1968221023713923
$ws.Range("K19").Value = 2.040694766216518
$ws.Range("L19").Value = 0.2165931591651571
$ws.Range("N19").Value = 1.390211291984201
$ws.Range("O19").Value = 4.7456234672226
$ws.Range("C20").Value = 0.2032217169418686
$ws.Range("D20").Value = 0.1539447817073381
$ws.Range("E20").Value = 0.1564856654987601
$ws.Range("F20").Value = 1.849789714709274
$ws.Range("G20").Value = 0.002489000905055993
$ws.Range("I20").Value = 1.299038219435054
$ws.Range("J20").Value = 0.197013499840935
$ws.Range("K20").Value = 2.093638862734281
$ws.Range("L20").Value = 0.2169435917018632
$ws.Range("N20").Value = 1.389973394596367
$ws.Range("O20").Value = 4.741330268832883
$ws.Range("C21").Value = 0.2081247564274378
$ws.Range("D21").Value = 0.1581504463494809
$ws.Range("E21").Value = 0.1583490023607759
$ws.Range("F21").Value = 1.850712281190312
$ws.Range("G21").Value = 0.0024856210104626
$ws.Range("I21").Value = 1.299156690492083
$ws.Range("J21").Value = 0.197746188235044
$ws.Range("K21").Value = 2.271703440832141
$ws.Range("L21").Value = 0.2182101604092921
$ws.Range("N21").Value = 1.389750741930939
$ws.Range("O21").Value = 4.729972802185387
$ws.Range("C22").Value = 0.2113869445672663
$ws.Range("D22").Value = 0.1609349793854449
$ws.Range("E22").Value = 0.1596139266521028
$ws.Range("F22").Value = 1.852157694874137
$ws.Range("G22").Value = 0.002483496405161128
$ws.Range("I22").Value = 1.299835000491676
$ws.Range("J22").Value = 0.1982916866974662
$ws.Range("K22").Value = 2.388200316929044
$ws.Range("L22").Value = 0.2191041632678647
$ws.Range("N22").Value = 1.390033676878076
$ws.Range("O22").Value = 4.724837680559745
$ws.Range("C23").Value = 0.2096405971395399
$ws.Range("D23").Value = 0.1594455648863118
$ws.Range("E23").Value = 0.1589345153376591
$ws.Range("F23").Value = 1.851309163327215
$ws.Range("G23").Value = 0.002484622650968411
$ws.Range("I23").Value = 1.299417985619087
$ws.Range("J23").Value = 0.197994455569571
$ws.Range("K23").Value = 2.326013438446012
$ws.Range("L23").Value = 0.2186209719425563
$ws.Range("N23").Value = 1.389843550359515
$ws.Range("O23").Value = 4.727368963612776
$ws.Range("C24").Value = 0.2031425611609734
$ws.Range("D24").Value = 0.1538766541831507
$ws.Range("E24").Value = 0.156456009071249
$ws.Range("F24").Value = 1.84978885800615
$ws.Range("G24").Value = 0.0024890575766401
$ws.Range("I24").Value = 1.299046425604026
$ws.Range("J24").Value = 0.1970026516735146
$ws.Range("K24").Value = 2.090730718025213
$ws.Range("L24").Value = 0.2169240116048243
$ws.Range("N24").Value = 1.38998428595913
$ws.Range("O24").Value = 4.74155449656692
$ws.Range("C25").Value = 0.1963843604388842
$ws.Range("D25").Value = 0.1480264645415161
$ws.Range("E25").Value = 0.1539863326862623
$ws.Range("F25").Value = 1.851765743500238
$ws.Range("G25").Value = 0.002494205785170451
$ws.Range("I25").Value = 1.301224671566111
$ws.Range("J25").Value = 0.1962199249687089
$ws.Range("K25").Value = 1.837541695693687
$ws.Range("L25").Value = 0.215379368589474
$ws.Range("N25").Value = 1.391984300616841
$ws.Range("O25").Value = 4.766666498949576
